$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the tab strip, named "Sheet3".
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Sheet3"

# Row 1: headers (shared with the other log sheets), plus the new "Date" column.
$ws.Range("A1").Value = "Solver Name"
$ws.Range("B1").Value = "Solver Variant"
$ws.Range("C1").Value = "Evaluation Function"
$ws.Range("D1").Value = "Move Gen Version"
$ws.Range("E1").Value = "Date"

# Row 2: values for this run, plus the new "???" placeholder under Date.
$ws.Range("A2").Value = "Minimax"
$ws.Range("B2").Value = "Basic"
$ws.Range("C2").Value = "ApplePieHeuristic"
$ws.Range("D2").Value = "v0.0"
$ws.Range("E2").Value = "???"

# Row 3 intentionally left blank (matches the other performance-log sheets).

# Row 4: column headers for the data table below.
$ws.Range("A4").Value = "Node Count"
$ws.Range("B4").Value = "Nodes Per Second"
$ws.Range("C4").Value = "Ply Count"
$ws.Range("D4").Value = "Plys Per Second"
$ws.Range("E4").Value = "Search Duration"

# Rows 5-10: measured data.
$data = @(
    @(206603, 158916.1875,     4, 3.07674503326416, 1.30007517337799),
    @(189903, 176677.578125,   4, 3.72142815589905, 1.07485604286194),
    @(190889, 187163.359375,   4, 3.92193078994751, 1.01990568637848),
    @(228744, 17920.630859375, 4, 0.313374429941177, 12.7642831802368),
    @(214282, 187889.625,      4, 3.50733399391174, 1.14046728610992),
    @(217141, 194081.03125,    4, 3.57520771026611, 1.11881613731384)
)

$row = 5
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("C$row").Value = $r[2]
    $ws.Range("D$row").Value = $r[3]
    $ws.Range("E$row").Value = $r[4]
    $row++
}
